$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "[59.82869319534586, 67.00213737464051]"
$ws.Range("T2").Value = "[47.08490829051013, 52.12261310729326]"
$ws.Range("L3").Value = "[58.40895593139782, 69.09465092375774]"
$ws.Range("T3").Value = "[46.93721128027376, 52.48049203289419]"
